$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master")

# --- Row 8 (EY_Pwd): update the password value and turn it into a hyperlink ---
$ws.Range("B8").Value = "P@ss.168"
$ws.Hyperlinks.Add($ws.Range("B8"), "https://etihad.okta.com/")
$ws.Range("B8").Style = "Hyperlink"

# --- Rows 14/15 (StartDate / EndDate): refresh the stored date text ---
$ws.Range("B14").Value = "16-SEP-2020"
$ws.Range("B15").Value = "16-SEP-2020"

# --- New rows 36-38: Power BI report URL settings ---
$ws.Range("A36").Value = "BIDAPReportURL"
$ws.Range("B36").Value = "https://app.powerbi.com/groups/137777c9-2dfd-406e-b75d-c8c3f6183177/reports/37f57a93-c504-411c-bfac-d39c71ffe917/ReportSection176e7ffa09d6b77bc221"
$ws.Range("C36").Value = "BI DAP Report URL"

$ws.Range("A37").Value = "BIGHAReportURL"
$ws.Range("B37").Value = "https://app.powerbi.com/groups/137777c9-2dfd-406e-b75d-c8c3f6183177/reports/37f57a93-c504-411c-bfac-d39c71ffe917/ReportSection176e7ffa09d6b77bc221"
$ws.Range("C37").Value = "BI GHA Repor tURL"

$ws.Range("A38").Value = "BIEvolutionReportURL"
$ws.Range("B38").Value = "https://app.powerbi.com/groups/137777c9-2dfd-406e-b75d-c8c3f6183177/reports/37f57a93-c504-411c-bfac-d39c71ffe917/ReportSectionbacdba451a1017289020?noSignUpCheck=1"
$ws.Range("C38").Value = "BI Evolution Report URL"

# --- Update the visible selection / scroll state on the Master sheet ---
$ws.Activate()
$ws.Range("B20").Select()
